$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "SME Performance Review EU" source block (rows 26-27)
# down to rows 32-33 by inserting 6 blank rows before row 26. This makes
# room for the new "Number of employees / Assets / Turnover" table that
# goes into rows 23-27.
$ws.Rows("26:31").Insert()

# New table header (row 23) - bold "title" style, matching the other
# table headers on this sheet (e.g. row 11, row 19).
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Font.Bold = $true

# New table body (rows 24-27) - plain "Normal" style.
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""

# Footer / source citation, now living at rows 32-33.
$ws.Range("A32").Value = "SME Performance Review EU"
$ws.Range("A32").Font.Bold = $true

$ws.Range("A33").Value = 'SME Performance Review EU, "SBA Fact sheet", 2013.  Available at http://ec.europa.eu/enterprise/policies/sme/facts-figures-analysis/performance-review/index_en.htm'
$ws.Range("A33").Font.Italic = $true
